$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "as of" date for this auto-update of stock data.
$newDate = "2026/01/15"

# Updated EBITDA readings per row (row 38 keeps its original EBITDA value).
$ebitda = @{
    2  = "8.07"
    8  = "8.86"
    14 = "3.23"
    20 = "13.76"
    26 = "11.92"
    32 = "29.24"
    44 = "16.25"
    50 = "12.39"
    56 = "32.47"
    62 = "11.95"
    68 = "13.27"
    74 = "18.73"
}

$rows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)

foreach ($r in $rows) {
    $aCell = $ws.Range("A$r")
    $aCell.Value = "'" + $newDate
    $aCell.Style = "Normal"

    if ($ebitda.ContainsKey($r)) {
        $bCell = $ws.Range("B$r")
        $bCell.Value = "'" + $ebitda[$r]
        $bCell.Style = "Normal"
    }
}
